$d = $word.ActiveDocument

# Helper: append one or more text chunks right after the first occurrence of
# $searchText within the whole document (used for the two table captions,
# "RAVDESS" and "MELD", which live in plain paragraphs, not table cells).
function AppendAfterText {
    param($searchText, [string[]]$parts)

    $found = $d.Content
    $found.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r = $d.Range($found.Start, $found.End)
    $r.Collapse(0)
    foreach ($p in $parts) {
        $r.InsertAfter($p)
        $r.Collapse(0)
    }
}

# Helper: append one or more text chunks right after the first occurrence of
# $searchText inside a specific table cell (table index, row, column are all
# 1-based, matching the Word object model).
function AppendAfterTextInCell {
    param($tableIdx, $row, $col, $searchText, [string[]]$parts)

    $table = $d.Tables.Item($tableIdx)
    $cellRange = $table.Cell($row, $col).Range
    $r = $d.Range($cellRange.Start, $cellRange.End)
    $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Collapse(0)
    foreach ($p in $parts) {
        $r.InsertAfter($p)
        $r.Collapse(0)
    }
}

# --- RAVDESS section -------------------------------------------------
AppendAfterText "RAVDESS" @(
    " (split train/test/val 0.7/0.2/0.1)",
    ", dla MFCC i eGEMAPS oraz dla wszystkich",
    " modeli poza CNN dane były standaryzowane i usuwane były atrybuty skorelowane"
)

AppendAfterTextInCell 1 2 1 "SVM" @(" (C=100)")
AppendAfterTextInCell 1 3 1 "RF" @(" (split po entropii, max-de", "pth 12)")
AppendAfterTextInCell 1 4 1 "LR" @(" (C=50, max-iter=1000)")
AppendAfterTextInCell 1 5 1 "MLP" @(" (max-iter=500)")
AppendAfterTextInCell 1 7 1 "GBT" @(" (subsample=0.5)")
AppendAfterTextInCell 1 8 1 "CNN" @(" (", "architektury w kodzie)")

# --- MELD section ------------------------------------------------------
AppendAfterText "MELD" @(
    " parametry modeli jak dla RAVDESS-a"
)

AppendAfterTextInCell 2 2 1 "SVM" @(" (class-weight=’balanced’)")
AppendAfterTextInCell 2 3 1 "RF" @(" ", "(class-weight=’balanced’)")
AppendAfterTextInCell 2 4 1 "LR" @(" ", "(class-weight=’balanced’)")
AppendAfterTextInCell 2 6 1 "DT" @(" ", "(class-weight=’balanced’)")
AppendAfterTextInCell 2 8 1 "CNN" @(
    " ",
    "(architektury w kodzie",
    ", dla MFCC zbalansowany zbiór poprzez resampling",
    ")"
)
